$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Court_Reports")

# Clear the example data rows (5 & 6) - content only, so formatting
# (e.g. the date-formatted Q column) is preserved.
$ws.Range("A5:Q6").ClearContents()

# Leave rows 5:6 selected, as recorded in the saved view state.
$ws.Rows("5:6").Select()
